$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105
$ws.Cells.Item(105, 1).Value = 7
$ws.Cells.Item(105, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(105, 3).Value = "Ñuble"
$ws.Cells.Item(105, 4).Value = 44567
$ws.Cells.Item(105, 5).Value = 16
$ws.Cells.Item(105, 6).Value = 100112032
$ws.Cells.Item(105, 7).Value = "Zapallo italiano"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 120
$ws.Cells.Item(105, 11).Value = 5000
$ws.Cells.Item(105, 12).Value = 5500
$ws.Cells.Item(105, 13).Value = 5250
$ws.Cells.Item(105, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(105, 15).Value = "Región del Maule"
$ws.Cells.Item(105, 16).Value = 66
$ws.Cells.Item(105, 17).Value = 80
$ws.Cells.Item(105, 18).Value = "Hortaliza"

# Row 106
$ws.Cells.Item(106, 1).Value = 7
$ws.Cells.Item(106, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(106, 3).Value = "Ñuble"
$ws.Cells.Item(106, 4).Value = 44567
$ws.Cells.Item(106, 5).Value = 16
$ws.Cells.Item(106, 6).Value = 100112032
$ws.Cells.Item(106, 7).Value = "Zapallo italiano"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Segunda"
$ws.Cells.Item(106, 10).Value = 120
$ws.Cells.Item(106, 11).Value = 4000
$ws.Cells.Item(106, 12).Value = 4500
$ws.Cells.Item(106, 13).Value = 4250
$ws.Cells.Item(106, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(106, 15).Value = "Región del Maule"
$ws.Cells.Item(106, 16).Value = 53
$ws.Cells.Item(106, 17).Value = 80
$ws.Cells.Item(106, 18).Value = "Hortaliza"

# Row 107
$ws.Cells.Item(107, 1).Value = 7
$ws.Cells.Item(107, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(107, 3).Value = "Ñuble"
$ws.Cells.Item(107, 4).Value = 44176
$ws.Cells.Item(107, 5).Value = 16
$ws.Cells.Item(107, 6).Value = 100112032
$ws.Cells.Item(107, 7).Value = "Zapallo italiano"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 120
$ws.Cells.Item(107, 11).Value = 7000
$ws.Cells.Item(107, 12).Value = 8000
$ws.Cells.Item(107, 13).Value = 7500
$ws.Cells.Item(107, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(107, 15).Value = "Región del Maule"
$ws.Cells.Item(107, 16).Value = 125
$ws.Cells.Item(107, 17).Value = 60
$ws.Cells.Item(107, 18).Value = "Hortaliza"

# Row 108
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 44370
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112032
$ws.Cells.Item(108, 7).Value = "Zapallo italiano"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 120
$ws.Cells.Item(108, 11).Value = 10000
$ws.Cells.Item(108, 12).Value = 11000
$ws.Cells.Item(108, 13).Value = 10500
$ws.Cells.Item(108, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(108, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(108, 16).Value = 175
$ws.Cells.Item(108, 17).Value = 60
$ws.Cells.Item(108, 18).Value = "Hortaliza"

# Row 109
$ws.Cells.Item(109, 1).Value = 7
$ws.Cells.Item(109, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(109, 3).Value = "Ñuble"
$ws.Cells.Item(109, 4).Value = 44475
$ws.Cells.Item(109, 5).Value = 16
$ws.Cells.Item(109, 6).Value = 100112032
$ws.Cells.Item(109, 7).Value = "Zapallo italiano"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 80
$ws.Cells.Item(109, 11).Value = 15000
$ws.Cells.Item(109, 12).Value = 16000
$ws.Cells.Item(109, 13).Value = 15500
$ws.Cells.Item(109, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(109, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(109, 16).Value = 310
$ws.Cells.Item(109, 17).Value = 50
$ws.Cells.Item(109, 18).Value = "Hortaliza"

# Row 110
$ws.Cells.Item(110, 1).Value = 7
$ws.Cells.Item(110, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(110, 3).Value = "Ñuble"
$ws.Cells.Item(110, 4).Value = 44508
$ws.Cells.Item(110, 5).Value = 16
$ws.Cells.Item(110, 6).Value = 100112032
$ws.Cells.Item(110, 7).Value = "Zapallo italiano"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 100
$ws.Cells.Item(110, 11).Value = 14000
$ws.Cells.Item(110, 12).Value = 15000
$ws.Cells.Item(110, 13).Value = 14500
$ws.Cells.Item(110, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(110, 15).Value = "Región del Maule"
$ws.Cells.Item(110, 16).Value = 242
$ws.Cells.Item(110, 17).Value = 60
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Row 111
$ws.Cells.Item(111, 1).Value = 7
$ws.Cells.Item(111, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111, 3).Value = "Ñuble"
$ws.Cells.Item(111, 4).Value = 44291
$ws.Cells.Item(111, 5).Value = 16
$ws.Cells.Item(111, 6).Value = 100112032
$ws.Cells.Item(111, 7).Value = "Zapallo italiano"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 120
$ws.Cells.Item(111, 11).Value = 6000
$ws.Cells.Item(111, 12).Value = 7000
$ws.Cells.Item(111, 13).Value = 6500
$ws.Cells.Item(111, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(111, 15).Value = "Región del Maule"
$ws.Cells.Item(111, 16).Value = 108
$ws.Cells.Item(111, 17).Value = 60
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Row 112
$ws.Cells.Item(112, 1).Value = 7
$ws.Cells.Item(112, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(112, 3).Value = "Ñuble"
$ws.Cells.Item(112, 4).Value = 44468
$ws.Cells.Item(112, 5).Value = 16
$ws.Cells.Item(112, 6).Value = 100112032
$ws.Cells.Item(112, 7).Value = "Zapallo italiano"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 120
$ws.Cells.Item(112, 11).Value = 14000
$ws.Cells.Item(112, 12).Value = 15000
$ws.Cells.Item(112, 13).Value = 14500
$ws.Cells.Item(112, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(112, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(112, 16).Value = 290
$ws.Cells.Item(112, 17).Value = 50
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Row 113
$ws.Cells.Item(113, 1).Value = 7
$ws.Cells.Item(113, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(113, 3).Value = "Ñuble"
$ws.Cells.Item(113, 4).Value = 44532
$ws.Cells.Item(113, 5).Value = 16
$ws.Cells.Item(113, 6).Value = 100112032
$ws.Cells.Item(113, 7).Value = "Zapallo italiano"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 100
$ws.Cells.Item(113, 11).Value = 8000
$ws.Cells.Item(113, 12).Value = 9000
$ws.Cells.Item(113, 13).Value = 8500
$ws.Cells.Item(113, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(113, 15).Value = "Región del Maule"
$ws.Cells.Item(113, 16).Value = 142
$ws.Cells.Item(113, 17).Value = 60
$ws.Cells.Item(113, 18).Value = "Hortaliza"

# Row 114
$ws.Cells.Item(114, 1).Value = 7
$ws.Cells.Item(114, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(114, 3).Value = "Ñuble"
$ws.Cells.Item(114, 4).Value = 44449
$ws.Cells.Item(114, 5).Value = 16
$ws.Cells.Item(114, 6).Value = 100112032
$ws.Cells.Item(114, 7).Value = "Zapallo italiano"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 100
$ws.Cells.Item(114, 11).Value = 16000
$ws.Cells.Item(114, 12).Value = 17000
$ws.Cells.Item(114, 13).Value = 16500
$ws.Cells.Item(114, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(114, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(114, 16).Value = 330
$ws.Cells.Item(114, 17).Value = 50
$ws.Cells.Item(114, 18).Value = "Hortaliza"

# Row 115
$ws.Cells.Item(115, 1).Value = 7
$ws.Cells.Item(115, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value = "Ñuble"
$ws.Cells.Item(115, 4).Value = 44210
$ws.Cells.Item(115, 5).Value = 16
$ws.Cells.Item(115, 6).Value = 100112032
$ws.Cells.Item(115, 7).Value = "Zapallo italiano"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 90
$ws.Cells.Item(115, 11).Value = 7000
$ws.Cells.Item(115, 12).Value = 8000
$ws.Cells.Item(115, 13).Value = 7333
$ws.Cells.Item(115, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(115, 15).Value = "Región del Maule"
$ws.Cells.Item(115, 16).Value = 122
$ws.Cells.Item(115, 17).Value = 60
$ws.Cells.Item(115, 18).Value = "Hortaliza"

# Row 116
$ws.Cells.Item(116, 1).Value = 7
$ws.Cells.Item(116, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(116, 3).Value = "Ñuble"
$ws.Cells.Item(116, 4).Value = 44526
$ws.Cells.Item(116, 5).Value = 16
$ws.Cells.Item(116, 6).Value = 100112032
$ws.Cells.Item(116, 7).Value = "Zapallo italiano"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 60
$ws.Cells.Item(116, 11).Value = 8000
$ws.Cells.Item(116, 12).Value = 9000
$ws.Cells.Item(116, 13).Value = 8500
$ws.Cells.Item(116, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(116, 15).Value = "Región del Maule"
$ws.Cells.Item(116, 16).Value = 142
$ws.Cells.Item(116, 17).Value = 60
$ws.Cells.Item(116, 18).Value = "Hortaliza"

# Row 117
$ws.Cells.Item(117, 1).Value = 7
$ws.Cells.Item(117, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(117, 3).Value = "Ñuble"
$ws.Cells.Item(117, 4).Value = 44271
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 100112032
$ws.Cells.Item(117, 7).Value = "Zapallo italiano"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 110
$ws.Cells.Item(117, 11).Value = 7000
$ws.Cells.Item(117, 12).Value = 7500
$ws.Cells.Item(117, 13).Value = 7295
$ws.Cells.Item(117, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(117, 15).Value = "Región del Maule"
$ws.Cells.Item(117, 16).Value = 122
$ws.Cells.Item(117, 17).Value = 60
$ws.Cells.Item(117, 18).Value = "Hortaliza"

# Row 118
$ws.Cells.Item(118, 1).Value = 7
$ws.Cells.Item(118, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(118, 3).Value = "Ñuble"
$ws.Cells.Item(118, 4).Value = 44425
$ws.Cells.Item(118, 5).Value = 16
$ws.Cells.Item(118, 6).Value = 100112032
$ws.Cells.Item(118, 7).Value = "Zapallo italiano"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 120
$ws.Cells.Item(118, 11).Value = 12000
$ws.Cells.Item(118, 12).Value = 13000
$ws.Cells.Item(118, 13).Value = 12500
$ws.Cells.Item(118, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(118, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(118, 16).Value = 250
$ws.Cells.Item(118, 17).Value = 50
$ws.Cells.Item(118, 18).Value = "Hortaliza"

# Row 119
$ws.Cells.Item(119, 1).Value = 7
$ws.Cells.Item(119, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(119, 3).Value = "Ñuble"
$ws.Cells.Item(119, 4).Value = 44218
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 6).Value = 100112032
$ws.Cells.Item(119, 7).Value = "Zapallo italiano"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 100
$ws.Cells.Item(119, 11).Value = 8000
$ws.Cells.Item(119, 12).Value = 9000
$ws.Cells.Item(119, 13).Value = 8650
$ws.Cells.Item(119, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(119, 15).Value = "Región del Maule"
$ws.Cells.Item(119, 16).Value = 144
$ws.Cells.Item(119, 17).Value = 60
$ws.Cells.Item(119, 18).Value = "Hortaliza"

# Row 120
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44250
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112032
$ws.Cells.Item(120, 7).Value = "Zapallo italiano"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 120
$ws.Cells.Item(120, 11).Value = 7000
$ws.Cells.Item(120, 12).Value = 8000
$ws.Cells.Item(120, 13).Value = 7500
$ws.Cells.Item(120, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(120, 15).Value = "Región del Maule"
$ws.Cells.Item(120, 16).Value = 125
$ws.Cells.Item(120, 17).Value = 60
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# Row 121
$ws.Cells.Item(121, 1).Value = 7
$ws.Cells.Item(121, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(121, 3).Value = "Ñuble"
$ws.Cells.Item(121, 4).Value = 44168
$ws.Cells.Item(121, 5).Value = 16
$ws.Cells.Item(121, 6).Value = 100112032
$ws.Cells.Item(121, 7).Value = "Zapallo italiano"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 160
$ws.Cells.Item(121, 11).Value = 6000
$ws.Cells.Item(121, 12).Value = 6500
$ws.Cells.Item(121, 13).Value = 6250
$ws.Cells.Item(121, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(121, 15).Value = "Región del Maule"
$ws.Cells.Item(121, 16).Value = 104
$ws.Cells.Item(121, 17).Value = 60
$ws.Cells.Item(121, 18).Value = "Hortaliza"

# Row 122
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 44447
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = 100112032
$ws.Cells.Item(122, 7).Value = "Zapallo italiano"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 160
$ws.Cells.Item(122, 11).Value = 16000
$ws.Cells.Item(122, 12).Value = 17000
$ws.Cells.Item(122, 13).Value = 16500
$ws.Cells.Item(122, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(122, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 16).Value = 330
$ws.Cells.Item(122, 17).Value = 50
$ws.Cells.Item(122, 18).Value = "Hortaliza"

# Row 123
$ws.Cells.Item(123, 1).Value = 7
$ws.Cells.Item(123, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(123, 3).Value = "Ñuble"
$ws.Cells.Item(123, 4).Value = 44553
$ws.Cells.Item(123, 5).Value = 16
$ws.Cells.Item(123, 6).Value = 100112032
$ws.Cells.Item(123, 7).Value = "Zapallo italiano"
$ws.Cells.Item(123, 8).Value = "Sin especificar"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 60
$ws.Cells.Item(123, 11).Value = 4000
$ws.Cells.Item(123, 12).Value = 4500
$ws.Cells.Item(123, 13).Value = 4250
$ws.Cells.Item(123, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(123, 15).Value = "Región del Maule"
$ws.Cells.Item(123, 16).Value = 71
$ws.Cells.Item(123, 17).Value = 60
$ws.Cells.Item(123, 18).Value = "Hortaliza"

# Row 124
$ws.Cells.Item(124, 1).Value = 7
$ws.Cells.Item(124, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(124, 3).Value = "Ñuble"
$ws.Cells.Item(124, 4).Value = 44167
$ws.Cells.Item(124, 5).Value = 16
$ws.Cells.Item(124, 6).Value = 100112032
$ws.Cells.Item(124, 7).Value = "Zapallo italiano"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 120
$ws.Cells.Item(124, 11).Value = 6000
$ws.Cells.Item(124, 12).Value = 6500
$ws.Cells.Item(124, 13).Value = 6250
$ws.Cells.Item(124, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(124, 15).Value = "Región del Maule"
$ws.Cells.Item(124, 16).Value = 104
$ws.Cells.Item(124, 17).Value = 60
$ws.Cells.Item(124, 18).Value = "Hortaliza"

# Row 125
$ws.Cells.Item(125, 1).Value = 7
$ws.Cells.Item(125, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125, 3).Value = "Ñuble"
$ws.Cells.Item(125, 4).Value = 44161
$ws.Cells.Item(125, 5).Value = 16
$ws.Cells.Item(125, 6).Value = 100112032
$ws.Cells.Item(125, 7).Value = "Zapallo italiano"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 130
$ws.Cells.Item(125, 11).Value = 6500
$ws.Cells.Item(125, 12).Value = 7000
$ws.Cells.Item(125, 13).Value = 6692
$ws.Cells.Item(125, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(125, 15).Value = "Región del Maule"
$ws.Cells.Item(125, 16).Value = 112
$ws.Cells.Item(125, 17).Value = 60
$ws.Cells.Item(125, 18).Value = "Hortaliza"

# Row 126
$ws.Cells.Item(126, 1).Value = 7
$ws.Cells.Item(126, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(126, 3).Value = "Ñuble"
$ws.Cells.Item(126, 4).Value = 44285
$ws.Cells.Item(126, 5).Value = 16
$ws.Cells.Item(126, 6).Value = 100112032
$ws.Cells.Item(126, 7).Value = "Zapallo italiano"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 120
$ws.Cells.Item(126, 11).Value = 8000
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = 8500
$ws.Cells.Item(126, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(126, 15).Value = "Región del Maule"
$ws.Cells.Item(126, 16).Value = 142
$ws.Cells.Item(126, 17).Value = 60
$ws.Cells.Item(126, 18).Value = "Hortaliza"

# Row 127
$ws.Cells.Item(127, 1).Value = 7
$ws.Cells.Item(127, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(127, 3).Value = "Ñuble"
$ws.Cells.Item(127, 4).Value = 44418
$ws.Cells.Item(127, 5).Value = 16
$ws.Cells.Item(127, 6).Value = 100112032
$ws.Cells.Item(127, 7).Value = "Zapallo italiano"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 120
$ws.Cells.Item(127, 11).Value = 11000
$ws.Cells.Item(127, 12).Value = 12000
$ws.Cells.Item(127, 13).Value = 11500
$ws.Cells.Item(127, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(127, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(127, 16).Value = 230
$ws.Cells.Item(127, 17).Value = 50
$ws.Cells.Item(127, 18).Value = "Hortaliza"

# Row 128
$ws.Cells.Item(128, 1).Value = 7
$ws.Cells.Item(128, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(128, 3).Value = "Ñuble"
$ws.Cells.Item(128, 4).Value = 44160
$ws.Cells.Item(128, 5).Value = 16
$ws.Cells.Item(128, 6).Value = 100112032
$ws.Cells.Item(128, 7).Value = "Zapallo italiano"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 190
$ws.Cells.Item(128, 11).Value = 6000
$ws.Cells.Item(128, 12).Value = 7000
$ws.Cells.Item(128, 13).Value = 6632
$ws.Cells.Item(128, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(128, 15).Value = "Región del Maule"
$ws.Cells.Item(128, 16).Value = 111
$ws.Cells.Item(128, 17).Value = 60
$ws.Cells.Item(128, 18).Value = "Hortaliza"

# Row 129
$ws.Cells.Item(129, 1).Value = 7
$ws.Cells.Item(129, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(129, 3).Value = "Ñuble"
$ws.Cells.Item(129, 4).Value = 44434
$ws.Cells.Item(129, 5).Value = 16
$ws.Cells.Item(129, 6).Value = 100112032
$ws.Cells.Item(129, 7).Value = "Zapallo italiano"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 160
$ws.Cells.Item(129, 11).Value = 14000
$ws.Cells.Item(129, 12).Value = 15000
$ws.Cells.Item(129, 13).Value = 14500
$ws.Cells.Item(129, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(129, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(129, 16).Value = 290
$ws.Cells.Item(129, 17).Value = 50
$ws.Cells.Item(129, 18).Value = "Hortaliza"

# Row 130
$ws.Cells.Item(130, 1).Value = 7
$ws.Cells.Item(130, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(130, 3).Value = "Ñuble"
$ws.Cells.Item(130, 4).Value = 44467
$ws.Cells.Item(130, 5).Value = 16
$ws.Cells.Item(130, 6).Value = 100112032
$ws.Cells.Item(130, 7).Value = "Zapallo italiano"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 120
$ws.Cells.Item(130, 11).Value = 11000
$ws.Cells.Item(130, 12).Value = 12000
$ws.Cells.Item(130, 13).Value = 11500
$ws.Cells.Item(130, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(130, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(130, 16).Value = 230
$ws.Cells.Item(130, 17).Value = 50
$ws.Cells.Item(130, 18).Value = "Hortaliza"

# Row 131
$ws.Cells.Item(131, 1).Value = 7
$ws.Cells.Item(131, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(131, 3).Value = "Ñuble"
$ws.Cells.Item(131, 4).Value = 44231
$ws.Cells.Item(131, 5).Value = 16
$ws.Cells.Item(131, 6).Value = 100112032
$ws.Cells.Item(131, 7).Value = "Zapallo italiano"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 120
$ws.Cells.Item(131, 11).Value = 8000
$ws.Cells.Item(131, 12).Value = 9000
$ws.Cells.Item(131, 13).Value = 8500
$ws.Cells.Item(131, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(131, 15).Value = "Región del Maule"
$ws.Cells.Item(131, 16).Value = 142
$ws.Cells.Item(131, 17).Value = 60
$ws.Cells.Item(131, 18).Value = "Hortaliza"

# Row 132
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(132, 3).Value = "Ñuble"
$ws.Cells.Item(132, 4).Value = 44259
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112032
$ws.Cells.Item(132, 7).Value = "Zapallo italiano"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 135
$ws.Cells.Item(132, 11).Value = 6500
$ws.Cells.Item(132, 12).Value = 7000
$ws.Cells.Item(132, 13).Value = 6722
$ws.Cells.Item(132, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(132, 15).Value = "Región del Maule"
$ws.Cells.Item(132, 16).Value = 112
$ws.Cells.Item(132, 17).Value = 60
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# Row 133
$ws.Cells.Item(133, 1).Value = 7
$ws.Cells.Item(133, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(133, 3).Value = "Ñuble"
$ws.Cells.Item(133, 4).Value = 44559
$ws.Cells.Item(133, 5).Value = 16
$ws.Cells.Item(133, 6).Value = 100112032
$ws.Cells.Item(133, 7).Value = "Zapallo italiano"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 100
$ws.Cells.Item(133, 11).Value = 5000
$ws.Cells.Item(133, 12).Value = 5500
$ws.Cells.Item(133, 13).Value = 5250
$ws.Cells.Item(133, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(133, 15).Value = "Región del Maule"
$ws.Cells.Item(133, 16).Value = 88
$ws.Cells.Item(133, 17).Value = 60
$ws.Cells.Item(133, 18).Value = "Hortaliza"

# Row 134
$ws.Cells.Item(134, 1).Value = 7
$ws.Cells.Item(134, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value = "Ñuble"
$ws.Cells.Item(134, 4).Value = 44341
$ws.Cells.Item(134, 5).Value = 16
$ws.Cells.Item(134, 6).Value = 100112032
$ws.Cells.Item(134, 7).Value = "Zapallo italiano"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 60
$ws.Cells.Item(134, 11).Value = 11000
$ws.Cells.Item(134, 12).Value = 12000
$ws.Cells.Item(134, 13).Value = 11500
$ws.Cells.Item(134, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(134, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(134, 16).Value = 192
$ws.Cells.Item(134, 17).Value = 60
$ws.Cells.Item(134, 18).Value = "Hortaliza"

# Row 135
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 44286
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112032
$ws.Cells.Item(135, 7).Value = "Zapallo italiano"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 120
$ws.Cells.Item(135, 11).Value = 8000
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 13).Value = 8500
$ws.Cells.Item(135, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(135, 15).Value = "Región del Maule"
$ws.Cells.Item(135, 16).Value = 142
$ws.Cells.Item(135, 17).Value = 60
$ws.Cells.Item(135, 18).Value = "Hortaliza"

# Row 136
$ws.Cells.Item(136, 1).Value = 7
$ws.Cells.Item(136, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(136, 3).Value = "Ñuble"
$ws.Cells.Item(136, 4).Value = 44208
$ws.Cells.Item(136, 5).Value = 16
$ws.Cells.Item(136, 6).Value = 100112032
$ws.Cells.Item(136, 7).Value = "Zapallo italiano"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 130
$ws.Cells.Item(136, 11).Value = 6500
$ws.Cells.Item(136, 12).Value = 7000
$ws.Cells.Item(136, 13).Value = 6731
$ws.Cells.Item(136, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(136, 15).Value = "Región del Maule"
$ws.Cells.Item(136, 16).Value = 112
$ws.Cells.Item(136, 17).Value = 60
$ws.Cells.Item(136, 18).Value = "Hortaliza"

# Row 137
$ws.Cells.Item(137, 1).Value = 7
$ws.Cells.Item(137, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(137, 3).Value = "Ñuble"
$ws.Cells.Item(137, 4).Value = 44264
$ws.Cells.Item(137, 5).Value = 16
$ws.Cells.Item(137, 6).Value = 100112032
$ws.Cells.Item(137, 7).Value = "Zapallo italiano"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 130
$ws.Cells.Item(137, 11).Value = 6500
$ws.Cells.Item(137, 12).Value = 7000
$ws.Cells.Item(137, 13).Value = 6769
$ws.Cells.Item(137, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(137, 15).Value = "Región del Maule"
$ws.Cells.Item(137, 16).Value = 113
$ws.Cells.Item(137, 17).Value = 60
$ws.Cells.Item(137, 18).Value = "Hortaliza"

# Row 138
$ws.Cells.Item(138, 1).Value = 7
$ws.Cells.Item(138, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(138, 3).Value = "Ñuble"
$ws.Cells.Item(138, 4).Value = 44322
$ws.Cells.Item(138, 5).Value = 16
$ws.Cells.Item(138, 6).Value = 100112032
$ws.Cells.Item(138, 7).Value = "Zapallo italiano"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 80
$ws.Cells.Item(138, 11).Value = 10000
$ws.Cells.Item(138, 12).Value = 11000
$ws.Cells.Item(138, 13).Value = 10500
$ws.Cells.Item(138, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(138, 15).Value = "Región del Maule"
$ws.Cells.Item(138, 16).Value = 175
$ws.Cells.Item(138, 17).Value = 60
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Row 139
$ws.Cells.Item(139, 1).Value = 7
$ws.Cells.Item(139, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(139, 3).Value = "Ñuble"
$ws.Cells.Item(139, 4).Value = 44391
$ws.Cells.Item(139, 5).Value = 16
$ws.Cells.Item(139, 6).Value = 100112032
$ws.Cells.Item(139, 7).Value = "Zapallo italiano"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 120
$ws.Cells.Item(139, 11).Value = 11000
$ws.Cells.Item(139, 12).Value = 12000
$ws.Cells.Item(139, 13).Value = 11500
$ws.Cells.Item(139, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(139, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(139, 16).Value = 192
$ws.Cells.Item(139, 17).Value = 60
$ws.Cells.Item(139, 18).Value = "Hortaliza"

# Row 140
$ws.Cells.Item(140, 1).Value = 7
$ws.Cells.Item(140, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(140, 3).Value = "Ñuble"
$ws.Cells.Item(140, 4).Value = 44396
$ws.Cells.Item(140, 5).Value = 16
$ws.Cells.Item(140, 6).Value = 100112032
$ws.Cells.Item(140, 7).Value = "Zapallo italiano"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 120
$ws.Cells.Item(140, 11).Value = 11000
$ws.Cells.Item(140, 12).Value = 12000
$ws.Cells.Item(140, 13).Value = 11500
$ws.Cells.Item(140, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(140, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(140, 16).Value = 230
$ws.Cells.Item(140, 17).Value = 50
$ws.Cells.Item(140, 18).Value = "Hortaliza"

# Row 141
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(141, 3).Value = "Ñuble"
$ws.Cells.Item(141, 4).Value = 44510
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 100112032
$ws.Cells.Item(141, 7).Value = "Zapallo italiano"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 120
$ws.Cells.Item(141, 11).Value = 10000
$ws.Cells.Item(141, 12).Value = 11000
$ws.Cells.Item(141, 13).Value = 10500
$ws.Cells.Item(141, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(141, 15).Value = "Región del Maule"
$ws.Cells.Item(141, 16).Value = 175
$ws.Cells.Item(141, 17).Value = 60
$ws.Cells.Item(141, 18).Value = "Hortaliza"

# Row 142
$ws.Cells.Item(142, 1).Value = 7
$ws.Cells.Item(142, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(142, 3).Value = "Ñuble"
$ws.Cells.Item(142, 4).Value = 44519
$ws.Cells.Item(142, 5).Value = 16
$ws.Cells.Item(142, 6).Value = 100112032
$ws.Cells.Item(142, 7).Value = "Zapallo italiano"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 100
$ws.Cells.Item(142, 11).Value = 8000
$ws.Cells.Item(142, 12).Value = 9000
$ws.Cells.Item(142, 13).Value = 8500
$ws.Cells.Item(142, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(142, 15).Value = "Región del Maule"
$ws.Cells.Item(142, 16).Value = 142
$ws.Cells.Item(142, 17).Value = 60
$ws.Cells.Item(142, 18).Value = "Hortaliza"

# Row 143
$ws.Cells.Item(143, 1).Value = 7
$ws.Cells.Item(143, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(143, 3).Value = "Ñuble"
$ws.Cells.Item(143, 4).Value = 44420
$ws.Cells.Item(143, 5).Value = 16
$ws.Cells.Item(143, 6).Value = 100112032
$ws.Cells.Item(143, 7).Value = "Zapallo italiano"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 120
$ws.Cells.Item(143, 11).Value = 11000
$ws.Cells.Item(143, 12).Value = 12000
$ws.Cells.Item(143, 13).Value = 11500
$ws.Cells.Item(143, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 230
$ws.Cells.Item(143, 17).Value = 50
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# Row 144
$ws.Cells.Item(144, 1).Value = 7
$ws.Cells.Item(144, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(144, 3).Value = "Ñuble"
$ws.Cells.Item(144, 4).Value = 44414
$ws.Cells.Item(144, 5).Value = 16
$ws.Cells.Item(144, 6).Value = 100112032
$ws.Cells.Item(144, 7).Value = "Zapallo italiano"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 120
$ws.Cells.Item(144, 11).Value = 9000
$ws.Cells.Item(144, 12).Value = 10000
$ws.Cells.Item(144, 13).Value = 9500
$ws.Cells.Item(144, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(144, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(144, 16).Value = 190
$ws.Cells.Item(144, 17).Value = 50
$ws.Cells.Item(144, 18).Value = "Hortaliza"

# Row 145
$ws.Cells.Item(145, 1).Value = 7
$ws.Cells.Item(145, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(145, 3).Value = "Ñuble"
$ws.Cells.Item(145, 4).Value = 44543
$ws.Cells.Item(145, 5).Value = 16
$ws.Cells.Item(145, 6).Value = 100112032
$ws.Cells.Item(145, 7).Value = "Zapallo italiano"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 200
$ws.Cells.Item(145, 11).Value = 6000
$ws.Cells.Item(145, 12).Value = 7000
$ws.Cells.Item(145, 13).Value = 6500
$ws.Cells.Item(145, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(145, 15).Value = "Región del Maule"
$ws.Cells.Item(145, 16).Value = 108
$ws.Cells.Item(145, 17).Value = 60
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# Row 146
$ws.Cells.Item(146, 1).Value = 7
$ws.Cells.Item(146, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(146, 3).Value = "Ñuble"
$ws.Cells.Item(146, 4).Value = 44278
$ws.Cells.Item(146, 5).Value = 16
$ws.Cells.Item(146, 6).Value = 100112032
$ws.Cells.Item(146, 7).Value = "Zapallo italiano"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 160
$ws.Cells.Item(146, 11).Value = 7000
$ws.Cells.Item(146, 12).Value = 8000
$ws.Cells.Item(146, 13).Value = 7500
$ws.Cells.Item(146, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(146, 15).Value = "Región del Maule"
$ws.Cells.Item(146, 16).Value = 125
$ws.Cells.Item(146, 17).Value = 60
$ws.Cells.Item(146, 18).Value = "Hortaliza"

# Row 147
$ws.Cells.Item(147, 1).Value = 7
$ws.Cells.Item(147, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(147, 3).Value = "Ñuble"
$ws.Cells.Item(147, 4).Value = 44187
$ws.Cells.Item(147, 5).Value = 16
$ws.Cells.Item(147, 6).Value = 100112032
$ws.Cells.Item(147, 7).Value = "Zapallo italiano"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 120
$ws.Cells.Item(147, 11).Value = 7500
$ws.Cells.Item(147, 12).Value = 8000
$ws.Cells.Item(147, 13).Value = 7750
$ws.Cells.Item(147, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(147, 15).Value = "Región del Maule"
$ws.Cells.Item(147, 16).Value = 129
$ws.Cells.Item(147, 17).Value = 60
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# Row 148
$ws.Cells.Item(148, 1).Value = 7
$ws.Cells.Item(148, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(148, 3).Value = "Ñuble"
$ws.Cells.Item(148, 4).Value = 44446
$ws.Cells.Item(148, 5).Value = 16
$ws.Cells.Item(148, 6).Value = 100112032
$ws.Cells.Item(148, 7).Value = "Zapallo italiano"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 160
$ws.Cells.Item(148, 11).Value = 16000
$ws.Cells.Item(148, 12).Value = 17000
$ws.Cells.Item(148, 13).Value = 16500
$ws.Cells.Item(148, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 330
$ws.Cells.Item(148, 17).Value = 50
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# Row 149
$ws.Cells.Item(149, 1).Value = 7
$ws.Cells.Item(149, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(149, 3).Value = "Ñuble"
$ws.Cells.Item(149, 4).Value = 44350
$ws.Cells.Item(149, 5).Value = 16
$ws.Cells.Item(149, 6).Value = 100112032
$ws.Cells.Item(149, 7).Value = "Zapallo italiano"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 120
$ws.Cells.Item(149, 11).Value = 9000
$ws.Cells.Item(149, 12).Value = 10000
$ws.Cells.Item(149, 13).Value = 9500
$ws.Cells.Item(149, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(149, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149, 16).Value = 158
$ws.Cells.Item(149, 17).Value = 60
$ws.Cells.Item(149, 18).Value = "Hortaliza"

# Row 150
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 44529
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = 100112032
$ws.Cells.Item(150, 7).Value = "Zapallo italiano"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 100
$ws.Cells.Item(150, 11).Value = 8000
$ws.Cells.Item(150, 12).Value = 9000
$ws.Cells.Item(150, 13).Value = 8500
$ws.Cells.Item(150, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(150, 15).Value = "Región del Maule"
$ws.Cells.Item(150, 16).Value = 142
$ws.Cells.Item(150, 17).Value = 60
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# Row 151
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(151, 3).Value = "Ñuble"
$ws.Cells.Item(151, 4).Value = 44413
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = 100112032
$ws.Cells.Item(151, 7).Value = "Zapallo italiano"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 160
$ws.Cells.Item(151, 11).Value = 8000
$ws.Cells.Item(151, 12).Value = 9000
$ws.Cells.Item(151, 13).Value = 8500
$ws.Cells.Item(151, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(151, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(151, 16).Value = 170
$ws.Cells.Item(151, 17).Value = 50
$ws.Cells.Item(151, 18).Value = "Hortaliza"

# Row 152
$ws.Cells.Item(152, 1).Value = 7
$ws.Cells.Item(152, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(152, 3).Value = "Ñuble"
$ws.Cells.Item(152, 4).Value = 44238
$ws.Cells.Item(152, 5).Value = 16
$ws.Cells.Item(152, 6).Value = 100112032
$ws.Cells.Item(152, 7).Value = "Zapallo italiano"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 120
$ws.Cells.Item(152, 11).Value = 8000
$ws.Cells.Item(152, 12).Value = 9000
$ws.Cells.Item(152, 13).Value = 8500
$ws.Cells.Item(152, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(152, 15).Value = "Región del Maule"
$ws.Cells.Item(152, 16).Value = 142
$ws.Cells.Item(152, 17).Value = 60
$ws.Cells.Item(152, 18).Value = "Hortaliza"

# Row 153
$ws.Cells.Item(153, 1).Value = 7
$ws.Cells.Item(153, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(153, 3).Value = "Ñuble"
$ws.Cells.Item(153, 4).Value = 44257
$ws.Cells.Item(153, 5).Value = 16
$ws.Cells.Item(153, 6).Value = 100112032
$ws.Cells.Item(153, 7).Value = "Zapallo italiano"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 110
$ws.Cells.Item(153, 11).Value = 5500
$ws.Cells.Item(153, 12).Value = 6000
$ws.Cells.Item(153, 13).Value = 5727
$ws.Cells.Item(153, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(153, 15).Value = "Región del Maule"
$ws.Cells.Item(153, 16).Value = 95
$ws.Cells.Item(153, 17).Value = 60
$ws.Cells.Item(153, 18).Value = "Hortaliza"

# Row 154
$ws.Cells.Item(154, 1).Value = 7
$ws.Cells.Item(154, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value = "Ñuble"
$ws.Cells.Item(154, 4).Value = 44411
$ws.Cells.Item(154, 5).Value = 16
$ws.Cells.Item(154, 6).Value = 100112032
$ws.Cells.Item(154, 7).Value = "Zapallo italiano"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 120
$ws.Cells.Item(154, 11).Value = 11000
$ws.Cells.Item(154, 12).Value = 12000
$ws.Cells.Item(154, 13).Value = 11500
$ws.Cells.Item(154, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(154, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(154, 16).Value = 230
$ws.Cells.Item(154, 17).Value = 50
$ws.Cells.Item(154, 18).Value = "Hortaliza"

# Row 155
$ws.Cells.Item(155, 1).Value = 7
$ws.Cells.Item(155, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(155, 3).Value = "Ñuble"
$ws.Cells.Item(155, 4).Value = 44175
$ws.Cells.Item(155, 5).Value = 16
$ws.Cells.Item(155, 6).Value = 100112032
$ws.Cells.Item(155, 7).Value = "Zapallo italiano"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 80
$ws.Cells.Item(155, 11).Value = 7500
$ws.Cells.Item(155, 12).Value = 8000
$ws.Cells.Item(155, 13).Value = 7750
$ws.Cells.Item(155, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(155, 15).Value = "Región del Maule"
$ws.Cells.Item(155, 16).Value = 129
$ws.Cells.Item(155, 17).Value = 60
$ws.Cells.Item(155, 18).Value = "Hortaliza"

# Row 156
$ws.Cells.Item(156, 1).Value = 7
$ws.Cells.Item(156, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value = "Ñuble"
$ws.Cells.Item(156, 4).Value = 44196
$ws.Cells.Item(156, 5).Value = 16
$ws.Cells.Item(156, 6).Value = 100112032
$ws.Cells.Item(156, 7).Value = "Zapallo italiano"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 120
$ws.Cells.Item(156, 11).Value = 6000
$ws.Cells.Item(156, 12).Value = 7000
$ws.Cells.Item(156, 13).Value = 6500
$ws.Cells.Item(156, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(156, 15).Value = "Región del Maule"
$ws.Cells.Item(156, 16).Value = 108
$ws.Cells.Item(156, 17).Value = 60
$ws.Cells.Item(156, 18).Value = "Hortaliza"

# Row 157
$ws.Cells.Item(157, 1).Value = 7
$ws.Cells.Item(157, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(157, 3).Value = "Ñuble"
$ws.Cells.Item(157, 4).Value = 44200
$ws.Cells.Item(157, 5).Value = 16
$ws.Cells.Item(157, 6).Value = 100112032
$ws.Cells.Item(157, 7).Value = "Zapallo italiano"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 120
$ws.Cells.Item(157, 11).Value = 6000
$ws.Cells.Item(157, 12).Value = 7000
$ws.Cells.Item(157, 13).Value = 6417
$ws.Cells.Item(157, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(157, 15).Value = "Región del Maule"
$ws.Cells.Item(157, 16).Value = 107
$ws.Cells.Item(157, 17).Value = 60
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# Row 158
$ws.Cells.Item(158, 1).Value = 7
$ws.Cells.Item(158, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(158, 3).Value = "Ñuble"
$ws.Cells.Item(158, 4).Value = 44188
$ws.Cells.Item(158, 5).Value = 16
$ws.Cells.Item(158, 6).Value = 100112032
$ws.Cells.Item(158, 7).Value = "Zapallo italiano"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 120
$ws.Cells.Item(158, 11).Value = 8000
$ws.Cells.Item(158, 12).Value = 9000
$ws.Cells.Item(158, 13).Value = 8500
$ws.Cells.Item(158, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(158, 15).Value = "Región del Maule"
$ws.Cells.Item(158, 16).Value = 142
$ws.Cells.Item(158, 17).Value = 60
$ws.Cells.Item(158, 18).Value = "Hortaliza"

# Row 159
$ws.Cells.Item(159, 1).Value = 7
$ws.Cells.Item(159, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(159, 3).Value = "Ñuble"
$ws.Cells.Item(159, 4).Value = 44258
$ws.Cells.Item(159, 5).Value = 16
$ws.Cells.Item(159, 6).Value = 100112032
$ws.Cells.Item(159, 7).Value = "Zapallo italiano"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 170
$ws.Cells.Item(159, 11).Value = 6000
$ws.Cells.Item(159, 12).Value = 7000
$ws.Cells.Item(159, 13).Value = 6529
$ws.Cells.Item(159, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(159, 15).Value = "Región del Maule"
$ws.Cells.Item(159, 16).Value = 109
$ws.Cells.Item(159, 17).Value = 60
$ws.Cells.Item(159, 18).Value = "Hortaliza"

# Row 160
$ws.Cells.Item(160, 1).Value = 7
$ws.Cells.Item(160, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(160, 3).Value = "Ñuble"
$ws.Cells.Item(160, 4).Value = 44298
$ws.Cells.Item(160, 5).Value = 16
$ws.Cells.Item(160, 6).Value = 100112032
$ws.Cells.Item(160, 7).Value = "Zapallo italiano"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 120
$ws.Cells.Item(160, 11).Value = 9000
$ws.Cells.Item(160, 12).Value = 10000
$ws.Cells.Item(160, 13).Value = 9500
$ws.Cells.Item(160, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 158
$ws.Cells.Item(160, 17).Value = 60
$ws.Cells.Item(160, 18).Value = "Hortaliza"

# Row 161
$ws.Cells.Item(161, 1).Value = 7
$ws.Cells.Item(161, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(161, 3).Value = "Ñuble"
$ws.Cells.Item(161, 4).Value = 44224
$ws.Cells.Item(161, 5).Value = 16
$ws.Cells.Item(161, 6).Value = 100112032
$ws.Cells.Item(161, 7).Value = "Zapallo italiano"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 130
$ws.Cells.Item(161, 11).Value = 7500
$ws.Cells.Item(161, 12).Value = 8000
$ws.Cells.Item(161, 13).Value = 7808
$ws.Cells.Item(161, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(161, 15).Value = "Región del Maule"
$ws.Cells.Item(161, 16).Value = 130
$ws.Cells.Item(161, 17).Value = 60
$ws.Cells.Item(161, 18).Value = "Hortaliza"

# Row 162
$ws.Cells.Item(162, 1).Value = 7
$ws.Cells.Item(162, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(162, 3).Value = "Ñuble"
$ws.Cells.Item(162, 4).Value = 44452
$ws.Cells.Item(162, 5).Value = 16
$ws.Cells.Item(162, 6).Value = 100112032
$ws.Cells.Item(162, 7).Value = "Zapallo italiano"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 160
$ws.Cells.Item(162, 11).Value = 16000
$ws.Cells.Item(162, 12).Value = 17000
$ws.Cells.Item(162, 13).Value = 16500
$ws.Cells.Item(162, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(162, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(162, 16).Value = 330
$ws.Cells.Item(162, 17).Value = 50
$ws.Cells.Item(162, 18).Value = "Hortaliza"

# Row 163
$ws.Cells.Item(163, 1).Value = 7
$ws.Cells.Item(163, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(163, 3).Value = "Ñuble"
$ws.Cells.Item(163, 4).Value = 44195
$ws.Cells.Item(163, 5).Value = 16
$ws.Cells.Item(163, 6).Value = 100112032
$ws.Cells.Item(163, 7).Value = "Zapallo italiano"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 120
$ws.Cells.Item(163, 11).Value = 6500
$ws.Cells.Item(163, 12).Value = 7000
$ws.Cells.Item(163, 13).Value = 6750
$ws.Cells.Item(163, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 112
$ws.Cells.Item(163, 17).Value = 60
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Row 164
$ws.Cells.Item(164, 1).Value = 7
$ws.Cells.Item(164, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(164, 3).Value = "Ñuble"
$ws.Cells.Item(164, 4).Value = 44536
$ws.Cells.Item(164, 5).Value = 16
$ws.Cells.Item(164, 6).Value = 100112032
$ws.Cells.Item(164, 7).Value = "Zapallo italiano"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 60
$ws.Cells.Item(164, 11).Value = 8000
$ws.Cells.Item(164, 12).Value = 9000
$ws.Cells.Item(164, 13).Value = 8500
$ws.Cells.Item(164, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 142
$ws.Cells.Item(164, 17).Value = 60
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# Row 165
$ws.Cells.Item(165, 1).Value = 7
$ws.Cells.Item(165, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(165, 3).Value = "Ñuble"
$ws.Cells.Item(165, 4).Value = 44511
$ws.Cells.Item(165, 5).Value = 16
$ws.Cells.Item(165, 6).Value = 100112032
$ws.Cells.Item(165, 7).Value = "Zapallo italiano"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 100
$ws.Cells.Item(165, 11).Value = 8000
$ws.Cells.Item(165, 12).Value = 9000
$ws.Cells.Item(165, 13).Value = 8500
$ws.Cells.Item(165, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 142
$ws.Cells.Item(165, 17).Value = 60
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Row 166
$ws.Cells.Item(166, 1).Value = 7
$ws.Cells.Item(166, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(166, 3).Value = "Ñuble"
$ws.Cells.Item(166, 4).Value = 44239
$ws.Cells.Item(166, 5).Value = 16
$ws.Cells.Item(166, 6).Value = 100112032
$ws.Cells.Item(166, 7).Value = "Zapallo italiano"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 60
$ws.Cells.Item(166, 11).Value = 8000
$ws.Cells.Item(166, 12).Value = 9000
$ws.Cells.Item(166, 13).Value = 8500
$ws.Cells.Item(166, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(166, 15).Value = "Región del Maule"
$ws.Cells.Item(166, 16).Value = 142
$ws.Cells.Item(166, 17).Value = 60
$ws.Cells.Item(166, 18).Value = "Hortaliza"

# Row 167
$ws.Cells.Item(167, 1).Value = 7
$ws.Cells.Item(167, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(167, 3).Value = "Ñuble"
$ws.Cells.Item(167, 4).Value = 44344
$ws.Cells.Item(167, 5).Value = 16
$ws.Cells.Item(167, 6).Value = 100112032
$ws.Cells.Item(167, 7).Value = "Zapallo italiano"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 120
$ws.Cells.Item(167, 11).Value = 10000
$ws.Cells.Item(167, 12).Value = 11000
$ws.Cells.Item(167, 13).Value = 10500
$ws.Cells.Item(167, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(167, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(167, 16).Value = 210
$ws.Cells.Item(167, 17).Value = 50
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# Row 168
$ws.Cells.Item(168, 1).Value = 7
$ws.Cells.Item(168, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(168, 3).Value = "Ñuble"
$ws.Cells.Item(168, 4).Value = 44568
$ws.Cells.Item(168, 5).Value = 16
$ws.Cells.Item(168, 6).Value = 100112032
$ws.Cells.Item(168, 7).Value = "Zapallo italiano"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 120
$ws.Cells.Item(168, 11).Value = 5000
$ws.Cells.Item(168, 12).Value = 5500
$ws.Cells.Item(168, 13).Value = 5250
$ws.Cells.Item(168, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 88
$ws.Cells.Item(168, 17).Value = 60
$ws.Cells.Item(168, 18).Value = "Hortaliza"

# Row 169
$ws.Cells.Item(169, 1).Value = 7
$ws.Cells.Item(169, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value = "Ñuble"
$ws.Cells.Item(169, 4).Value = 44568
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = 100112032
$ws.Cells.Item(169, 7).Value = "Zapallo italiano"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Segunda"
$ws.Cells.Item(169, 10).Value = 60
$ws.Cells.Item(169, 11).Value = 4000
$ws.Cells.Item(169, 12).Value = 4500
$ws.Cells.Item(169, 13).Value = 4250
$ws.Cells.Item(169, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(169, 15).Value = "Región del Maule"
$ws.Cells.Item(169, 16).Value = 53
$ws.Cells.Item(169, 17).Value = 80
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Row 170
$ws.Cells.Item(170, 1).Value = 7
$ws.Cells.Item(170, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(170, 3).Value = "Ñuble"
$ws.Cells.Item(170, 4).Value = 44463
$ws.Cells.Item(170, 5).Value = 16
$ws.Cells.Item(170, 6).Value = 100112032
$ws.Cells.Item(170, 7).Value = "Zapallo italiano"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 120
$ws.Cells.Item(170, 11).Value = 14000
$ws.Cells.Item(170, 12).Value = 15000
$ws.Cells.Item(170, 13).Value = 14500
$ws.Cells.Item(170, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(170, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(170, 16).Value = 290
$ws.Cells.Item(170, 17).Value = 50
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# Row 171
$ws.Cells.Item(171, 1).Value = 7
$ws.Cells.Item(171, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(171, 3).Value = "Ñuble"
$ws.Cells.Item(171, 4).Value = 44365
$ws.Cells.Item(171, 5).Value = 16
$ws.Cells.Item(171, 6).Value = 100112032
$ws.Cells.Item(171, 7).Value = "Zapallo italiano"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 120
$ws.Cells.Item(171, 11).Value = 9000
$ws.Cells.Item(171, 12).Value = 10000
$ws.Cells.Item(171, 13).Value = 9500
$ws.Cells.Item(171, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(171, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(171, 16).Value = 158
$ws.Cells.Item(171, 17).Value = 60
$ws.Cells.Item(171, 18).Value = "Hortaliza"

# Row 172
$ws.Cells.Item(172, 1).Value = 7
$ws.Cells.Item(172, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(172, 3).Value = "Ñuble"
$ws.Cells.Item(172, 4).Value = 44454
$ws.Cells.Item(172, 5).Value = 16
$ws.Cells.Item(172, 6).Value = 100112032
$ws.Cells.Item(172, 7).Value = "Zapallo italiano"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 160
$ws.Cells.Item(172, 11).Value = 16000
$ws.Cells.Item(172, 12).Value = 17000
$ws.Cells.Item(172, 13).Value = 16500
$ws.Cells.Item(172, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(172, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(172, 16).Value = 330
$ws.Cells.Item(172, 17).Value = 50
$ws.Cells.Item(172, 18).Value = "Hortaliza"

# Row 173
$ws.Cells.Item(173, 1).Value = 7
$ws.Cells.Item(173, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(173, 3).Value = "Ñuble"
$ws.Cells.Item(173, 4).Value = 44194
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = 100112032
$ws.Cells.Item(173, 7).Value = "Zapallo italiano"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 60
$ws.Cells.Item(173, 11).Value = 6500
$ws.Cells.Item(173, 12).Value = 7000
$ws.Cells.Item(173, 13).Value = 6750
$ws.Cells.Item(173, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(173, 15).Value = "Región del Maule"
$ws.Cells.Item(173, 16).Value = 112
$ws.Cells.Item(173, 17).Value = 60
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# Row 174
$ws.Cells.Item(174, 1).Value = 7
$ws.Cells.Item(174, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(174, 3).Value = "Ñuble"
$ws.Cells.Item(174, 4).Value = 44518
$ws.Cells.Item(174, 5).Value = 16
$ws.Cells.Item(174, 6).Value = 100112032
$ws.Cells.Item(174, 7).Value = "Zapallo italiano"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 100
$ws.Cells.Item(174, 11).Value = 8000
$ws.Cells.Item(174, 12).Value = 9000
$ws.Cells.Item(174, 13).Value = 8500
$ws.Cells.Item(174, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(174, 15).Value = "Región del Maule"
$ws.Cells.Item(174, 16).Value = 142
$ws.Cells.Item(174, 17).Value = 60
$ws.Cells.Item(174, 18).Value = "Hortaliza"

# Row 175
$ws.Cells.Item(175, 1).Value = 7
$ws.Cells.Item(175, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(175, 3).Value = "Ñuble"
$ws.Cells.Item(175, 4).Value = 44272
$ws.Cells.Item(175, 5).Value = 16
$ws.Cells.Item(175, 6).Value = 100112032
$ws.Cells.Item(175, 7).Value = "Zapallo italiano"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 150
$ws.Cells.Item(175, 11).Value = 7000
$ws.Cells.Item(175, 12).Value = 7500
$ws.Cells.Item(175, 13).Value = 7267
$ws.Cells.Item(175, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(175, 15).Value = "Región del Maule"
$ws.Cells.Item(175, 16).Value = 121
$ws.Cells.Item(175, 17).Value = 60
$ws.Cells.Item(175, 18).Value = "Hortaliza"

$ws.Range("D105:D175").NumberFormat = "YYYY-MM-DD HH:MM:SS"
